$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(76, 1).Value = "2025-10-04T12:38:26.997476"
$ws.Cells.Item(76, 2).Value = "TRADING_ATTEMPT"
$ws.Cells.Item(76, 3).Value = "SOL"
$ws.Cells.Item(76, 4).Value = "UNKNOWN"
$ws.Cells.Item(76, 5).Value = 233.064580246805
$ws.Cells.Item(76, 11).Value = "ATTEMPT"
$ws.Cells.Item(76, 12).Value = "Attempting trade 1/7"

$ws.Cells.Item(77, 1).Value = "2025-10-04T12:38:28.975269"
$ws.Cells.Item(77, 2).Value = "POSITION_FAILED"
$ws.Cells.Item(77, 3).Value = "SOL"
$ws.Cells.Item(77, 4).Value = "UNKNOWN"
$ws.Cells.Item(77, 11).Value = "FAILED"
$ws.Cells.Item(77, 12).Value = "Trade execution failed for trade 1"

$ws.Cells.Item(78, 1).Value = "2025-10-04T12:38:29.008107"
$ws.Cells.Item(78, 2).Value = "TRADING_ATTEMPT"
$ws.Cells.Item(78, 3).Value = "BTC"
$ws.Cells.Item(78, 4).Value = "UNKNOWN"
$ws.Cells.Item(78, 5).Value = 122250.151868428
$ws.Cells.Item(78, 11).Value = "ATTEMPT"
$ws.Cells.Item(78, 12).Value = "Attempting trade 2/7"

$ws.Cells.Item(79, 1).Value = "2025-10-04T12:38:30.745560"
$ws.Cells.Item(79, 2).Value = "POSITION_FAILED"
$ws.Cells.Item(79, 3).Value = "BTC"
$ws.Cells.Item(79, 4).Value = "UNKNOWN"
$ws.Cells.Item(79, 11).Value = "FAILED"
$ws.Cells.Item(79, 12).Value = "Trade execution failed for trade 2"

$ws.Cells.Item(80, 1).Value = "2025-10-04T12:38:30.778527"
$ws.Cells.Item(80, 2).Value = "TRADING_ATTEMPT"
$ws.Cells.Item(80, 3).Value = "ETH"
$ws.Cells.Item(80, 4).Value = "UNKNOWN"
$ws.Cells.Item(80, 5).Value = 4515.759068159723
$ws.Cells.Item(80, 11).Value = "ATTEMPT"
$ws.Cells.Item(80, 12).Value = "Attempting trade 3/7"

$ws.Cells.Item(81, 1).Value = "2025-10-04T12:38:32.474841"
$ws.Cells.Item(81, 2).Value = "POSITION_FAILED"
$ws.Cells.Item(81, 3).Value = "ETH"
$ws.Cells.Item(81, 4).Value = "UNKNOWN"
$ws.Cells.Item(81, 11).Value = "FAILED"
$ws.Cells.Item(81, 12).Value = "Trade execution failed for trade 3"

$ws.Cells.Item(82, 1).Value = "2025-10-04T12:38:32.507311"
$ws.Cells.Item(82, 2).Value = "TRADING_ATTEMPT"
$ws.Cells.Item(82, 3).Value = "ARB"
$ws.Cells.Item(82, 4).Value = "UNKNOWN"
$ws.Cells.Item(82, 5).Value = 0.4505143051833715
$ws.Cells.Item(82, 11).Value = "ATTEMPT"
$ws.Cells.Item(82, 12).Value = "Attempting trade 4/7"

$ws.Cells.Item(83, 1).Value = "2025-10-04T12:38:34.377397"
$ws.Cells.Item(83, 2).Value = "POSITION_FAILED"
$ws.Cells.Item(83, 3).Value = "ARB"
$ws.Cells.Item(83, 4).Value = "UNKNOWN"
$ws.Cells.Item(83, 11).Value = "FAILED"
$ws.Cells.Item(83, 12).Value = "Trade execution failed for trade 4"

$ws.Cells.Item(84, 1).Value = "2025-10-04T12:38:34.412320"
$ws.Cells.Item(84, 2).Value = "TRADING_ATTEMPT"
$ws.Cells.Item(84, 3).Value = "SUI"
$ws.Cells.Item(84, 4).Value = "UNKNOWN"
$ws.Cells.Item(84, 5).Value = 3.582136289260382
$ws.Cells.Item(84, 11).Value = "ATTEMPT"
$ws.Cells.Item(84, 12).Value = "Attempting trade 5/7"

$ws.Cells.Item(85, 1).Value = "2025-10-04T12:38:36.296444"
$ws.Cells.Item(85, 2).Value = "POSITION_FAILED"
$ws.Cells.Item(85, 3).Value = "SUI"
$ws.Cells.Item(85, 4).Value = "UNKNOWN"
$ws.Cells.Item(85, 11).Value = "FAILED"
$ws.Cells.Item(85, 12).Value = "Trade execution failed for trade 5"

$ws.Cells.Item(86, 1).Value = "2025-10-04T12:38:36.331327"
$ws.Cells.Item(86, 2).Value = "TRADING_ATTEMPT"
$ws.Cells.Item(86, 3).Value = "XRP"
$ws.Cells.Item(86, 4).Value = "UNKNOWN"
$ws.Cells.Item(86, 5).Value = 3.041366390190445
$ws.Cells.Item(86, 11).Value = "ATTEMPT"
$ws.Cells.Item(86, 12).Value = "Attempting trade 6/7"

$ws.Cells.Item(87, 1).Value = "2025-10-04T12:38:38.148810"
$ws.Cells.Item(87, 2).Value = "POSITION_FAILED"
$ws.Cells.Item(87, 3).Value = "XRP"
$ws.Cells.Item(87, 4).Value = "UNKNOWN"
$ws.Cells.Item(87, 11).Value = "FAILED"
$ws.Cells.Item(87, 12).Value = "Trade execution failed for trade 6"

$ws.Cells.Item(88, 1).Value = "2025-10-04T12:38:38.183243"
$ws.Cells.Item(88, 2).Value = "TRADING_ATTEMPT"
$ws.Cells.Item(88, 3).Value = "AAVE"
$ws.Cells.Item(88, 4).Value = "UNKNOWN"
$ws.Cells.Item(88, 5).Value = 290.9949382388142
$ws.Cells.Item(88, 11).Value = "ATTEMPT"
$ws.Cells.Item(88, 12).Value = "Attempting trade 7/7"

$ws.Cells.Item(89, 1).Value = "2025-10-04T12:38:40.071041"
$ws.Cells.Item(89, 2).Value = "POSITION_FAILED"
$ws.Cells.Item(89, 3).Value = "AAVE"
$ws.Cells.Item(89, 4).Value = "UNKNOWN"
$ws.Cells.Item(89, 11).Value = "FAILED"
$ws.Cells.Item(89, 12).Value = "Trade execution failed for trade 7"

"done"
